# Apply the "experiment with internetarchive extractor" cleanup:
#  - strip trailing "." / ". :" / " :" punctuation left over from MARC-style
#    fields (title, subtitle, place, publisher) across the existing rows
#  - blank out the placeholder "[s.n.] :" / "," / "n.d.." publication cells
#    on row 6 (no real publisher data available for that record)
#  - append five new catalogue rows (12-16) pulled in by the extractor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($row, $col, $value) {
    # Plain text cell write.
    $ws.Cells.Item($row, $col).Value2 = $value
}

function Set-IdText($row, $col, $value) {
    # Purely-numeric IDs/years need a leading apostrophe so Excel keeps
    # them as text, matching every other "numeric-looking" field already
    # stored as a string on this sheet (not a real number).
    $ws.Cells.Item($row, $col).Value2 = "'" + $value
}

function Set-Blank($row, $col) {
    # Write an empty string, then touch a no-op format property so the
    # engine keeps the cell as a present-but-empty record instead of
    # dropping it from the row entirely (and without allocating a new
    # cell style in the process).
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value2 = ""
    $cell.Font.Bold = $false
}

# --- Row 2: Abercrombie, Inquiries concerning ... ---
Set-Text 2 3 "Inquiries concerning the intellectual powers and the investigation of truth"
Set-Text 2 5 "London"
Set-Text 2 6 "John Murray"

# --- Row 3: Darwin supplementary copy ---
Set-Text 3 3 "Inquiries concerning the intellectual powers and the investigation of truth"
Set-Text 3 4 "[Supplementary material in Charles Darwin's copy]"

# --- Row 4: Agassiz, Lake Superior ---
Set-Text 4 3 "Lake Superior: its character, vegetation, and animals, compared with those of other similar regions"
Set-Text 4 5 "Boston"
Set-Text 4 6 "Gould, Kendall & Lincoln"

# --- Row 5: Darwin supplementary copy ---
Set-Text 5 3 "Lake Superior: its character, vegetation, and animals, compared with those of other similar regions"
Set-Text 5 4 "[Supplementary material in Charles Darwin's copy]"

# --- Row 6: Agassiz, Contributions ... (no real pub info) ---
Set-Text 6 3 "Contributions to the natural history of the United States of North America"
Set-Blank 6 5
Set-Blank 6 6
Set-Blank 6 7

# --- Row 7: Darwin supplementary copy ---
Set-Text 7 3 "Contributions to the natural history of the United States of North America"
Set-Text 7 4 "[Supplementary material in Charles Darwin's copy]"

# --- Row 8: Barker-Webb, Histoire naturelle ---
Set-Text 8 3 "Histoire naturelle des Îles Canaries"
Set-Text 8 5 "Paris"
Set-Text 8 6 "Béthune"

# --- Row 9: Darwin supplementary copy ---
Set-Text 9 3 "Histoire naturelle des Îles Canaries"
Set-Text 9 4 "[Supplementary material in Charles Darwin's copy]"

# --- Row 10: Bechstein, Naturgeschichte der Stubenvögel ---
Set-Text 10 3 "Naturgeschichte der Stubenvögel"
Set-Text 10 5 "Halle"
Set-Text 10 6 "Hennemann"

# --- Row 11: Darwin supplementary copy ---
Set-Text 11 3 "Naturgeschichte der Stubenvögel"
Set-Text 11 4 "[Supplementary material in Charles Darwin's copy]"

# --- Row 12: Bell, The anatomy and philosophy of expression ---
Set-IdText 12 1 "1000241"
Set-Text   12 2 "Bell, Charles"
Set-Text   12 3 "The anatomy and philosophy of expression"
Set-Blank  12 4
Set-Text   12 5 "London"
Set-Text   12 6 "John Murray"
Set-IdText 12 7 "1844"
Set-Text   12 8 "Signature E. Darwin 1844 to Ch. Darwin Nov. 28 1866.; Location: Cambridge.; Identifier: Bell1844nu34M.; Public number: 0085."
Set-Blank  12 9

# --- Row 13: Darwin supplementary copy ---
Set-IdText 13 1 "1000242"
Set-Text   13 2 "Darwin, Charles"
Set-Text   13 3 "The anatomy and philosophy of expression"
Set-Text   13 4 "[Supplementary material in Charles Darwin's copy]"
Set-Blank  13 5
Set-Blank  13 6
Set-Blank  13 7
Set-Text   13 8 "Identifier: Bell1844nu34M_MS."
Set-Blank  13 9

# --- Row 14: Bernhardi, über den Begriff der Pflanzenart ---
Set-IdText 14 1 "1000116"
Set-Text   14 2 "Bernhardi, Johann Jacob"
Set-Text   14 3 "über den Begriff der Pflanzenart und seine Anwendung"
Set-Blank  14 4
Set-Text   14 5 "Erfurt"
Set-Text   14 6 "Friedrich Wilhelm Otto"
Set-IdText 14 7 "1834"
Set-Text   14 8 "Location: Cambridge.; Identifier: Bernhardi1834ez86I.; Public number: 0098."
Set-Blank  14 9

# --- Row 15: Darwin supplementary copy ---
Set-IdText 15 1 "1000117"
Set-Text   15 2 "Darwin, Charles"
Set-Text   15 3 "über den Begriff der Pflanzenart und seine Anwendung"
Set-Text   15 4 "[Supplementary material in Charles Darwin's copy]"
Set-Blank  15 5
Set-Blank  15 6
Set-Blank  15 7
Set-Text   15 8 "Identifier: Bernhardi1834ez86I_MS."
Set-Blank  15 9

# --- Row 16: Boitard, Les Pigeons de volière et de colombier ---
Set-IdText 16 1 "1000141"
Set-Text   16 2 "Boitard, Pierre"
Set-Text   16 3 "Les Pigeons de volière et de colombier, ou Histoire naturelle des pigeons domestiques"
Set-Blank  16 4
Set-Text   16 5 "Paris"
Set-Text   16 6 "Audot & Corbié"
Set-IdText 16 7 "1824"
Set-Text   16 8 "Pre-Beagle.; Location: Cambridge.; Identifier: Boitard1824jp50S.; Public number: 0118."
Set-Blank  16 9
